$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was published for this market/variety; it sorts
# to the top of the data block (row 523), pushing the existing rows
# 523-551 down by one (to 524-552).
$ws.Rows.Item(523).Insert()

$ws.Cells.Item(523, 1).Value = 9
$ws.Cells.Item(523, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(523, 3).Value = "Metropolitana"
$ws.Cells.Item(523, 4).Value = 45008
$ws.Cells.Item(523, 5).Value = 13
$ws.Cells.Item(523, 6).Value = 100112032
$ws.Cells.Item(523, 7).Value = "Zapallo italiano"
$ws.Cells.Item(523, 8).Value = "Sin especificar"
$ws.Cells.Item(523, 9).Value = "Primera"
$ws.Cells.Item(523, 10).Value = 430
$ws.Cells.Item(523, 11).Value = 6000
$ws.Cells.Item(523, 12).Value = 7000
$ws.Cells.Item(523, 13).Value = 6500
$ws.Cells.Item(523, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(523, 15).Value = "Región Metropolitana"
$ws.Cells.Item(523, 16).Value = 130
$ws.Cells.Item(523, 17).Value = 50
$ws.Cells.Item(523, 18).Value = "Hortaliza"
